# Anonymize "fedcore" -> "approach" and tidy up the header borders on the
# merged title row of both sheets (quality_comparison / computational_comparison).
#
# Border edge indices used below match real Excel constants:
#   7  = xlEdgeLeft
#   8  = xlEdgeTop
#   9  = xlEdgeBottom
#   10 = xlEdgeRight
# LineStyle 1 = xlContinuous (thin line), -4142 = xlLineStyleNone

$wb = $excel.ActiveWorkbook

function Set-HeaderCellBorders {
    param($cell, $hasRight)
    $cell.ClearFormats()
    if ($hasRight) {
        $cell.Borders.Item(10).LineStyle = 1   # right
    }
    $cell.Borders.Item(8).LineStyle = 1        # top
    $cell.Borders.Item(9).LineStyle = 1        # bottom
}

# ---- Sheet: quality_comparison ----
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-HeaderCellBorders $ws1.Range("C1") $false
Set-HeaderCellBorders $ws1.Range("D1") $true

$ws1.Range("C2").Value = "approach"

# ---- Sheet: computational_comparison ----
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-HeaderCellBorders $ws2.Range("C1") $false
Set-HeaderCellBorders $ws2.Range("D1") $true
Set-HeaderCellBorders $ws2.Range("F1") $false
Set-HeaderCellBorders $ws2.Range("G1") $true

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell.
$ws2.Range("G5").ClearContents()
